# Generate Report for Handoff
# Update the "Latest Handoff Date/Datetime" values for the row corresponding
# to file 4999b8bf-edcd-4b0b-bbf6-e5582e2da185.md in each worksheet, since a
# newer handoff run completed for that file.

$wb = $excel.ActiveWorkbook

# Overview sheet: column D = "Latest Handoff Date" for row 6.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D6").Value = "2016-03-23 22:40:47"

# zh-cn sheet: column E = "Latest Handoff Datetime" for row 6.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E6").Value = "2016-03-23 22:40:43"

# de-de sheet: column E = "Latest Handoff Datetime" for row 6.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E6").Value = "2016-03-23 22:40:47"
